# The calibration data in the "Summary" sheet (time (s), AA1/AA2/AA3
# curvature) needs to be sorted chronologically (ascending by column A,
# "time (s)"), as a result of (re-)performing the needle calibration.
# Only the data rows (2-8) are reordered; the header row (1) is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:D8")
$sortKey = $ws.Range("A2:A8")

$dataRange.Sort($sortKey, 1)
